# Add probability rows for true angle 1, 2, 8, and 9
#
# The table models a discrete probability distribution (spread 0.1/0.2/0.4/0.2/0.1)
# of "measured angle" given a "true angle" (rows 11-19, true angles 1-9, columns C-K
# = measured angle 1-9). Rows for true angle 3-7 were already filled in. This fills
# in the edge rows (1, 2, 8, 9) where the bell-shaped spread gets clipped against the
# edge of the table, so the clipped probability mass piles up on the boundary cell.
# The boundary ("edge") cell of each row is additionally highlighted with a distinct
# fill color (reddish/pink) to flag it, matching the existing orange highlight style
# already used for the interior spread cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BGR integer for RGB FFC000 (existing orange highlight fill already used by the sheet)
$orange = 49407
# BGR integer for RGB FF7C80 (new pink/red highlight fill for the clipped edge cell)
$pink = 8420607

function Set-ProbCell($addr, $value, $color) {
    $cell = $ws.Range($addr)
    $cell.Value2 = $value
    $cell.Interior.Color = $color
}

# True angle = 1 (row 11): spread clipped on the left, mass piles on C11
Set-ProbCell "C11" 0.7 $pink
Set-ProbCell "D11" 0.2 $orange
Set-ProbCell "E11" 0.1 $orange

# True angle = 2 (row 12): spread clipped on the left, mass piles on C12
Set-ProbCell "C12" 0.3 $pink
Set-ProbCell "D12" 0.39999999999999991 $orange
Set-ProbCell "E12" 0.2 $orange
Set-ProbCell "F12" 0.1 $orange

# True angle = 8 (row 18): spread clipped on the right, mass piles on K18
Set-ProbCell "H18" 0.1 $orange
Set-ProbCell "I18" 0.2 $orange
Set-ProbCell "J18" 0.39999999999999991 $orange
Set-ProbCell "K18" 0.3 $pink

# True angle = 9 (row 19): spread clipped on the right, mass piles on K19
Set-ProbCell "I19" 0.1 $orange
Set-ProbCell "J19" 0.2 $orange
Set-ProbCell "K19" 0.7 $pink

# The author's selection moved off B9 when saving; clear the lingering selection
# marker left over from editing by reselecting the sheet's natural home cell.
$ws.Range("A1").Select()
